$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 42 becomes the former row-43 values (Homo sapiens / Human)
$ws.Range("A42").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B42").Value = "Homo sapiens"
$ws.Range("C42").Value = "Human"
$ws.Range("D42").Value = "Human"

# Row 43 becomes the former row-42 values (unassigned)
$ws.Range("A43").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B43").Value = "unassigned"
$ws.Range("C43").Value = "unassigned"
$ws.Range("D43").Value = "unassigned"

# Row 55 becomes the former row-56 values (unassigned)
$ws.Range("A55").Value = "5e733a21f67e541f28ed4bf4fe025044"
$ws.Range("B55").Value = "unassigned"
$ws.Range("C55").Value = "unassigned"
$ws.Range("D55").Value = "unassigned"

# Row 56 becomes the former row-55 values (Centropristis striata / Black sea bass / Teleost Fish)
$ws.Range("A56").Value = "975b1dbdc7405f6e27bf63893e91e0ed"
$ws.Range("B56").Value = "Centropristis striata"
$ws.Range("C56").Value = "Black sea bass"
$ws.Range("D56").Value = "Teleost Fish"
